$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B33 changes from a text "3" to a true number 3
$ws.Range("B33").Value = 3

# New row 34
$ws.Range("A34").Value = "Ruilin"

# B34 stays a text string "3" (not numeric) per the source diff
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "3"

$ws.Range("C34").Value = "无"
$ws.Range("D34").Value = "DFT"
$ws.Range("E34").Value = "MET"
$ws.Range("F34").Value = "2e6daeb6-f5b1-42e4-9927-e16202e5fb2e"
$ws.Range("G34").Value = "H1cWzoxA-_annotated.xlsx"
$ws.Range("H34").Value = 'For example ,when I use the cr dataset, "python sc_main.py --network_type exp_context_fusion --context_fusion_method wblock --model_dir_suffix training --dataset_type cr --gpu 0 " the result is not the 84.48 as the paper,I could only get 84.30 after several times.'
